# Insert a new weekly record at row 8, pushing the existing rows (old 8..29)
# down to (9..30). The new row replicates the constant columns for this
# product/market and carries the newest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8:29 down by one to make room for the new weekly record.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the latest week's data.
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44575
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 100114007
$ws.Range("G8").Value = "Jengibre"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 14500
$ws.Range("N8").Value = "`$/caja 13 kilos"
$ws.Range("O8").Value = "Perú"
$ws.Range("P8").Value = 1115
$ws.Range("Q8").Value = 13
$ws.Range("R8").Value = "Hortaliza"
